{"js": "// The document contains several short paragraphs like \"<id>p060v_1</id>\"\n// that were split across multiple runs (e.g. a run for \"<id>\", a separate\n// run for the value \"p060v_1\" (sometimes itself split further, e.g.\n// \"p060v_\" + \"4\"), and a run for \"</id>\"). The edit collapses each of\n// these into a single run (using the formatting of the \"<id>\" run) that\n// contains the full concatenated text, e.g. \"<id>p060v_1</id>\".\n\nconst body = context.document.body;\n\n// Find every paragraph that starts an \"<id>...</id>\" marker.\nconst markers = body.search(\"<id>\", { matchCase: true });\nmarkers.load(\"text\");\nawait context.sync();\n\n// Collect the full (logical) text of each matching paragraph first -\n// search results reference fixed ranges, so we must gather the target\n// text up front before mutating the document.\nconst fullTexts = [];\nfor (let i = 0; i < markers.items.length; i++) {\n  const paragraph = markers.items[i].paragraphs.getFirst();\n  paragraph.load(\"text\");\n  await context.sync();\n  fullTexts.push(paragraph.text);\n}\n\n// For each paragraph, re-locate it by its exact full text and replace the\n// (multi-run) match in place with the same text. Word/Office.js merges a\n// \"Replace\" insertion into a single run that takes on the formatting of\n// the first run it overwrote, which is exactly the \"<id>\" run's styling\n// (Courier New / color 7f6000 / sz 18).\nfor (const fullText of fullTexts) {\n  const hits = body.search(fullText, { matchCase: true });\n  hits.load(\"text\");\n  await context.sync();\n\n  if (hits.items.length > 0) {\n    hits.items[0].insertText(fullText, \"Replace\");\n    await context.sync();\n  }\n}\n", "ps1": "# The document contains several short paragraphs like \"<id>p060v_1</id>\"\n# whose text is split across multiple runs (e.g. one run for the literal\n# \"<id>\", a separate run for the value \"p060v_1\" - sometimes itself split\n# further, e.g. \"p060v_\" + \"4\" - and a run for the literal \"</id>\"). This\n# collapses each such paragraph down to a single run (keeping the\n# formatting of the \"<id>\" run, i.e. Courier New / color 7f6000 / sz 18)\n# that contains the full concatenated text, e.g. \"<id>p060v_1</id>\".\n\n$d = $word.ActiveDocument\n\n# Locate every \"<id>\" marker in the document (in document order).\n$content = $d.Content\n$find = $content.Find\n$find.ClearFormatting()\n$find.Text = \"<id>\"\n$find.Forward = $true\n$find.Wrap = 0  # wdFindStop - do not wrap back to the start\n\n$starts = @()\nwhile ($find.Execute()) {\n    $starts += $content.Start\n    [void]$content.Collapse(0)  # wdCollapseEnd - continue searching after this hit\n}\n\nforeach ($s in $starts) {\n    $para = $d.Range($s, $s).Paragraphs(1)\n    $prange = $para.Range\n    [void]$prange.MoveEnd(1, -1)  # wdCharacter - drop the trailing paragraph mark\n\n    # First 4 characters are always the literal \"<id>\" marker run.\n    $markerRange = $d.Range($prange.Start, $prange.Start + 4)\n    # Everything else in the paragraph is the id value plus the closing \"</id>\".\n    $remainderRange = $d.Range($prange.Start + 4, $prange.End)\n    $remainderText = $remainderRange.Text\n\n    # Delete the extra runs, then append their text onto the \"<id>\" run so\n    # the whole marker collapses into a single run using its formatting.\n    $remainderRange.Delete()\n    [void]$markerRange.Collapse(0)  # wdCollapseEnd\n    $markerRange.InsertAfter($remainderText)\n}\n"}
